# ------------------------------------------------------------------
# [ADDITIONAL SCRAPING] split the single "ODI Batting" sheet workbook
# into three sheets:
#   1. "Player Info"        - basic player bio fields
#   2. "ODI Batting"         - existing per-match batting log, but the
#                              MATCH_CARD_LINK column is replaced with a
#                              bare MATCH_CODE number
#   3. "ODI Batting Extra"   - additional per-match batting detail
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# The workbook currently has exactly one sheet: "ODI Batting".
$wsBatting = $wb.Worksheets.Item(1)

# --- 1. Insert the new "Player Info" sheet *before* "ODI Batting" ----
$wsInfo = $wb.Worksheets.Add($wsBatting)
$wsInfo.Name = "Player Info"

$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $infoHeaders.Length; $col++) {
    $cell = $wsInfo.Cells.Item(1, $col)
    $cell.Value = $infoHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$infoRow = @("3728", "Mohammad Shahzad", "Right Handed", "Does Not Bowl | Unknown")
$wsInfo.Range("A2:D2").NumberFormat = "@"
for ($col = 1; $col -le $infoRow.Length; $col++) {
    $wsInfo.Cells.Item(2, $col).Value = $infoRow[$col - 1]
}

# --- 2. Insert the new "ODI Batting Extra" sheet *after* "ODI Batting" ----
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsExtra = $wb.Worksheets.Add($null, $wsBatting)
$wsExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $extraHeaders.Length; $col++) {
    $cell = $wsExtra.Cells.Item(1, $col)
    $cell.Value = $extraHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# MATCH_CODE (A) and MAN_OF_MATCH (F) are always text; the middle
# columns are text too (even though most values look numeric), except
# BATTING_POSITION (B) which is a genuine number when present.
$wsExtra.Range("A2:A21").NumberFormat = "@"
$wsExtra.Range("C2:F21").NumberFormat = "@"

# Blank cells in the source data are genuine empty strings (not
# missing/null cells) - use "" as the sentinel for "empty" and $null
# only for the BATTING_POSITION numeric column when a row has no
# batting-position entry at all.
$extraRows = @(
    @("4145", 1,     "3",  "0", "12.37%", "NO"),
    @("4154", $null, "",   "",  "",       "NO"),
    @("4160", 1,     "0",  "0", "1.69%",  "NO"),
    @("4163", 1,     "6",  "2", "25.35%", "YES"),
    @("4164", 1,     "11", "2", "40.78%", "YES"),
    @("4190", 1,     "0",  "0", "",       "NO"),
    @("4192", 1,     "0",  "0", "0.79%",  "NO"),
    @("4195", 1,     "4",  "1", "13.65%", "NO"),
    @("4198", 1,     "4",  "0", "14.51%", "NO"),
    @("4200", 1,     "2",  "0", "7.78%",  "NO"),
    @("4202", 1,     "8",  "0", "21.54%", "NO"),
    @("4203", $null, "",   "",  "",       "NO"),
    @("4257", 1,     "6",  "0", "26.06%", "NO"),
    @("4262", $null, "",   "",  "",       "NO"),
    @("4267", 1,     "0",  "1", "2.78%",  "NO"),
    @("4290", 1,     "7",  "1", "20.45%", "NO"),
    @("4299", 1,     "0",  "0", "1.45%",  "NO"),
    @("4301", 1,     "16", "0", "33.11%", "NO"),
    @("4306", 1,     "0",  "0", "",       "NO"),
    @("4309", 1,     "1",  "0", "4.61%",  "NO")
)

$r = 2
foreach ($row in $extraRows) {
    for ($col = 1; $col -le 6; $col++) {
        $v = $row[$col - 1]
        if ($null -ne $v) {
            $wsExtra.Cells.Item($r, $col).Value = $v
        }
    }
    $r++
}

# --- 3. Update the "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ----
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"

$lastRow = $wsBatting.UsedRange.Rows.Count
$wsBatting.Range("D2:D" + $lastRow).NumberFormat = "@"
for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $wsBatting.Cells.Item($i, 4)
    $val = $cell.Value2
    if ($val -match "MatchCode=(\d+)") {
        $cell.Value = $matches[1]
    }
}

# --- 4. Restore the originally active sheet/tab ----
$wsInfo2 = $wb.Worksheets.Item("Player Info")
$wsInfo2.Activate()
